# ---------------------------------------------------------------------------
# "formatted mean avg deviations table in excel"
#
#   * "Deviations" is rebuilt as a formatted "Table 1": a title, an italic
#     subtitle, trial headers (Trial 1..Trial 10 + Average of 10 Trials),
#     and the per-scenario rows grouped under "PI(t)D(t)" / "Traditional
#     PID" row-group labels.
#   * The old "Scenario Max" column is pulled out into a brand new
#     worksheet "Max Errors" (placed right after "Deviations"), using the
#     same row-group layout.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$devSheet = $wb.Worksheets.Item("Deviations")

# ---------------------------------------------------------------------------
# 1. Snapshot the existing numbers before anything is overwritten.
#    Rows 2-7 (col A = scenario index 0..5): B..L = 11 trial values,
#    M = the scenario's max deviation.
# ---------------------------------------------------------------------------
$trialVals = @{}
$maxVals = @{}
for ($r = 2; $r -le 7; $r++) {
    $rowVals = @()
    for ($c = 2; $c -le 12; $c++) {
        $rowVals += $devSheet.Cells.Item($r, $c).Value2
    }
    $trialVals[$r] = $rowVals
    $maxVals[$r] = $devSheet.Cells.Item($r, 13).Value2
}

# ---------------------------------------------------------------------------
# 2. Insert the new "Max Errors" sheet right after "Deviations".
# ---------------------------------------------------------------------------
$maxSheet = $wb.Worksheets.Add($null, $devSheet)
$maxSheet.Name = "Max Errors"

# ---------------------------------------------------------------------------
# 3. Clear out the old "Deviations" contents so the new layout can be
#    written from scratch.
# ---------------------------------------------------------------------------
$devSheet.Cells.Clear()

$scenarioLabels = @("Scenario 1", "Scenario 2", "Scenario 3")
$shortLabels    = @("S1", "S2", "S3")
$groupLabels    = @("PI(t)D(t)", "Traditional PID")

# ---------------------------------------------------------------------------
# 4. Write all the text/number content first (values only, no styling) so
#    later formatting passes can be applied on whole (unioned) ranges.
# ---------------------------------------------------------------------------

# -- "Deviations" / "Table 1" --------------------------------------------
$devSheet.Range("A1").Value = "Table 1"
$devSheet.Range("A2").Value = "Mean Average Deviation From Target Path (in)"

$trialHeaders = @("Trial 1","Trial 2","Trial 3","Trial 4","Trial 5","Trial 6","Trial 7","Trial 8","Trial 9","Trial 10")
for ($i = 0; $i -lt 10; $i++) {
    $devSheet.Cells.Item(3, 3 + $i).Value = $trialHeaders[$i]
}
$devSheet.Cells.Item(3, 13).Value = "Average of 10 Trials"

$devSheet.Range("A4").Value = "PI(t)D(t)"
$devSheet.Range("A7").Value = "Traditional PID"

$srcRows = @(2, 3, 4, 5, 6, 7)
$dstRows = @(4, 5, 6, 7, 8, 9)
for ($i = 0; $i -lt 6; $i++) {
    $dstRow = $dstRows[$i]
    $srcRow = $srcRows[$i]
    $label = $scenarioLabels[$i % 3]

    $devSheet.Cells.Item($dstRow, 2).Value = $label

    $rowVals = $trialVals[$srcRow]
    for ($c = 0; $c -lt 10; $c++) {
        $devSheet.Cells.Item($dstRow, 3 + $c).Value = $rowVals[$c]
    }
    $devSheet.Cells.Item($dstRow, 13).Value = $rowVals[10]
}

# -- "Max Errors" ----------------------------------------------------------
$maxSheet.Range("C1").Value = "Max."
$maxDstRows = @(2, 3, 4, 5, 6, 7)
for ($i = 0; $i -lt 6; $i++) {
    $dstRow = $maxDstRows[$i]
    $srcRow = $srcRows[$i]
    $maxSheet.Cells.Item($dstRow, 2).Value = $shortLabels[$i % 3]
    $maxSheet.Cells.Item($dstRow, 3).Value = $maxVals[$srcRow]
}
$maxSheet.Range("A2").Value = "PI(t)D(t)"
$maxSheet.Range("A5").Value = "Traditional PID"

# ---------------------------------------------------------------------------
# 5. Formatting. Build each distinct font once (Bold first keeps the
#    font table smallest), applied across a Union of every range that
#    needs it, then layer on alignment / number format / merges.
# ---------------------------------------------------------------------------

# Bold Arial 10: title, trial headers, row-group labels, scenario labels.
$boldRange = $excel.Union($devSheet.Range("A1:B1"), $devSheet.Range("C3:M3"))
$boldRange = $excel.Union($boldRange, $devSheet.Range("A4:A9"))
$boldRange = $excel.Union($boldRange, $devSheet.Range("B4:B9"))
$boldRange = $excel.Union($boldRange, $maxSheet.Range("C1"))
$boldRange = $excel.Union($boldRange, $maxSheet.Range("A2:A7"))
$boldRange = $excel.Union($boldRange, $maxSheet.Range("B2:B7"))
$boldRange.Font.Bold = $true
$boldRange.Font.Name = "Arial"
$boldRange.Font.Size = 10

# Italic Arial 10: subtitle row.
$devSheet.Range("A2:N2").Font.Italic = $true
$devSheet.Range("A2:N2").Font.Name = "Arial"
$devSheet.Range("A2:N2").Font.Size = 10

# Regular Arial 10: blank corner cells + numeric data cells.
$regRange = $excel.Union($devSheet.Range("A3:B3"), $devSheet.Range("C4:M9"))
$regRange = $excel.Union($regRange, $maxSheet.Range("C2:C7"))
$regRange.Font.Name = "Arial"
$regRange.Font.Size = 10

# -- number formats ---------------------------------------------------------
$devSheet.Range("C4:M9").NumberFormat = "0.000"
$maxSheet.Range("C2:C7").NumberFormat = "0.000"

# -- alignment ---------------------------------------------------------------
$devSheet.Range("C3:M3").HorizontalAlignment = -4108
$devSheet.Range("C3:M3").VerticalAlignment = -4160

$devSheet.Range("A4:A9").HorizontalAlignment = -4131
$devSheet.Range("B4:B9").HorizontalAlignment = -4131
$devSheet.Range("A7:A9").WrapText = $true

$devSheet.Range("C4:M9").HorizontalAlignment = -4108

$devSheet.Range("A2:N2").HorizontalAlignment = -4131

$maxSheet.Range("A2:A7").HorizontalAlignment = -4131
$maxSheet.Range("B2:B7").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 6. Merges.
# ---------------------------------------------------------------------------
[void] $devSheet.Range("A2:N2").Merge()
[void] $devSheet.Range("A4:A6").Merge()
[void] $devSheet.Range("A7:A9").Merge()

[void] $maxSheet.Range("A2:A4").Merge()
[void] $maxSheet.Range("A5:A7").Merge()

# ---------------------------------------------------------------------------
# 7. Column widths, selection, misc.
# ---------------------------------------------------------------------------
$devSheet.Columns.Item(1).ColumnWidth = 12.29
$devSheet.Columns.Item(2).ColumnWidth = 10
$devSheet.Columns.Item(13).ColumnWidth = 18.53
$devSheet.Columns.Item(14).ColumnWidth = 12.59

$maxSheet.Columns.Item(2).ColumnWidth = 3.53

$devSheet.Range("A1").Select()
